# Update absenteeism data rows 2-11 (Colaborador records) to new values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A=Colaborador_id, B=Colaborador_nome, C=Departamento,
#          D=Motivo_da_ausência, E=Horas_de_ausência, F=Data_da_ausência, G=Salário
$rows = @(
    @{ Row=2;  A=9616;  B="Ana Clara da Mata";     C="Recursos Humanos";       D="Consulta médica";    E=8; F=45096; G=4293.36 },
    @{ Row=3;  A=9089;  B="Sr. Breno da Cruz";      C="Jurídico";              D="Consulta médica";    E=5; F=45101; G=6215.63 },
    @{ Row=4;  A=8688;  B="João Vitor da Cruz";     C="Atendimento ao Cliente";D="Problemas pessoais";  E=5; F=45081; G=3428.81 },
    @{ Row=5;  A=66936; B="João Vitor Moreira";     C="Recursos Humanos";      D="Doença";              E=2; F=45097; G=8613.11 },
    @{ Row=6;  A=30967; B="Laura Costela";          C="Engenharia";            D="Outros";              E=5; F=45085; G=8971.34 },
    @{ Row=7;  A=26806; B="Nina Ramos";             C="TI";                    D="Doença";              E=5; F=45086; G=6677.24 },
    @{ Row=8;  A=96308; B="Rafael Araújo";          C="P&D";                   D="Doença";              E=3; F=45089; G=12273.34 },
    @{ Row=9;  A=28965; B="Eduardo Ferreira";       C="Financeiro";            D="Doença";              E=5; F=45090; G=11810.61 },
    @{ Row=10; A=84225; B="Daniela Lima";           C="Recursos Humanos";      D="Viagem de negócios";  E=1; F=45106; G=7894.42 },
    @{ Row=11; A=27844; B="Luiz Henrique Rocha";    C="Recursos Humanos";      D="Doença";              E=3; F=45079; G=10670.1 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row
    $ws.Cells.Item($rowNum, 1).Value = $r.A
    $ws.Cells.Item($rowNum, 2).Value = $r.B
    $ws.Cells.Item($rowNum, 3).Value = $r.C
    $ws.Cells.Item($rowNum, 4).Value = $r.D
    $ws.Cells.Item($rowNum, 5).Value = $r.E
    $ws.Cells.Item($rowNum, 6).Value = $r.F
    $ws.Cells.Item($rowNum, 7).Value = $r.G
}
